$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mentor names
$ws.Range("A2").Value = "Mentor1 Surname"
$ws.Range("A3").Value = "Mentor2 Surname"

# Update availability hours
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 1

# Clear old mentor rows (A4:B14) - remove values but keep formatting
$ws.Range("A4:B14").ClearContents()

# Extend formatting down through rows 15-17 in column A (copy style from A14/A13)
$ws.Range("A13:A14").Copy()
$ws.Range("A15:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Ensure B column cleared for new rows too (should already be empty)
$ws.Range("B15:B17").ClearContents()

# Update selection to match target diff (B3 selected)
$ws.Range("B3").Select()
